$wb = $excel.ActiveWorkbook

function Set-TextValue($range, $text) {
    # Force the cell to be stored as Text (matching the source workbook,
    # which stores these numeric-looking strings as shared strings, not
    # numbers), then strip the resulting quote-prefix style so the cell
    # keeps the default (General) style, same as the rest of the sheet.
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

# --- Sheet: Restricciones_del_follower (index 3) ---
$ws = $wb.Worksheets.Item(3)

Set-TextValue $ws.Range("A2") "0.6691555555555564 - 0.45600000000000007y_1 + 0.7857777777777776y_2"
Set-TextValue $ws.Range("B2") "-0.6691555555555564"
Set-TextValue $ws.Range("D2") "0.35"
Set-TextValue $ws.Range("E2") "0"
Set-TextValue $ws.Range("F2") "1.1"

Set-TextValue $ws.Range("A3") "-2.4206000000000003 + 0.6759999999999999y_1 - 0.4679999999999999y_2"
Set-TextValue $ws.Range("B3") "-1.5793999999999997"
Set-TextValue $ws.Range("D3") "0.91"
Set-TextValue $ws.Range("E3") "9.7"
Set-TextValue $ws.Range("F3") "0"

Set-TextValue $ws.Range("A4") "20.482599999999998 - 2x - 1.1959999999999997y_1 + 0.8280000000000007y_2"
Set-TextValue $ws.Range("B4") "-36.4826"
Set-TextValue $ws.Range("D4") "0.78"
Set-TextValue $ws.Range("E4") "6.3"
Set-TextValue $ws.Range("F4") "0"

Set-TextValue $ws.Range("A5") "-69.3236 + 8x + 0.856y_1 - 0.208y_2"
Set-TextValue $ws.Range("B5") "20.653599999999997"
Set-TextValue $ws.Range("D5") "0.57"
Set-TextValue $ws.Range("E5") "0"
Set-TextValue $ws.Range("F5") "0.8999999999999999"

Set-TextValue $ws.Range("A6") "5.187911111111112 - 2x - 1.072y_1 + 1.3404444444444443y_2"
Set-TextValue $ws.Range("B6") "-6.812088888888888"
Set-TextValue $ws.Range("D6") "0.0"
Set-TextValue $ws.Range("E6") "0"
Set-TextValue $ws.Range("F6") "1.4000000000000001"

# --- Sheet: Punto_modificado (index 4) ---
$ws = $wb.Worksheets.Item(4)
Set-TextValue $ws.Range("A2") "8.1"
Set-TextValue $ws.Range("B2") "5.0"
Set-TextValue $ws.Range("C2") "2.05"

# --- Sheet: Vector_bf (index 5) ---
$ws = $wb.Worksheets.Item(5)
Set-TextValue $ws.Range("A2") "0.9894"
Set-TextValue $ws.Range("A3") "-0.3764222222222228"

# --- Sheet: Vector_BF (index 6) ---
$ws = $wb.Worksheets.Item(6)
Set-TextValue $ws.Range("A2") "13.6"
Set-TextValue $ws.Range("A3") "3.977599999999999"
Set-TextValue $ws.Range("A4") "-2.6768000000000054"

# --- Sheet: Vector_Alpha (index 7) ---
# These two cells are genuine numbers in both the original and the edited
# workbook, so set them as real numeric values (no text coercion needed).
$ws = $wb.Worksheets.Item(7)
$ws.Range("A2").Value = 1.62
$ws.Range("A3").Value = 2.34
